$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reformat D2:D51 as text first, and reset style back to Normal afterward so that
# numeric-looking price strings (e.g. "1.00", "34.90") are preserved verbatim as text
# instead of being silently coerced to numbers (which would drop trailing zeros /
# introduce float rounding) - while keeping the cell style index unchanged (matches
# the target workbook, which carries no explicit style on these cells).
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = '69.368.34'
$ws.Range("E2").Value = '  -1.01%  '

$ws.Range("D3").Value = '3.533.15'
$ws.Range("E3").Value = '  -2.16%  '

$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.20%  '

$ws.Range("D5").Value = '196.59'
$ws.Range("E5").Value = '  +0.53%  '

$ws.Range("D6").Value = '585.49'
$ws.Range("E6").Value = '  -3.34%  '

$ws.Range("D7").Value = '0.613'
$ws.Range("E7").Value = '  -2.09%  '

$ws.Range("D9").Value = '0.206'
$ws.Range("E9").Value = '  -0.40%  '

$ws.Range("D10").Value = '0.628'
$ws.Range("E10").Value = '  -3.75%  '

$ws.Range("D11").Value = '51.88'
$ws.Range("E11").Value = '  -3.86%  '

$ws.Range("D12").Value = '0.0000288'
$ws.Range("E12").Value = '  -5.73%  '

$ws.Range("E13").Value = '  -2.50%  '

$ws.Range("D14").Value = '684.76'
$ws.Range("E14").Value = '  +15.37%  '

$ws.Range("D15").Value = '4.090.95'
$ws.Range("E15").Value = '  -2.15%  '

$ws.Range("D16").Value = '69.451.09'
$ws.Range("E16").Value = '  -1.16%  '

$ws.Range("B17").Value = 'Uniswap'
$ws.Range("C17").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D17").Value = '12.44'
$ws.Range("E17").Value = '  -5.52%  '

$ws.Range("B18").Value = 'Chainlink'
$ws.Range("C18").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D18").Value = '18.57'
$ws.Range("E18").Value = '  -3.84%  '

$ws.Range("B19").Value = 'WrappedEther'
$ws.Range("C19").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D19").Value = '3.510.47'
$ws.Range("E19").Value = '  -2.69%  '

$ws.Range("D20").Value = '0.121'
$ws.Range("E20").Value = '  -0.86%  '

$ws.Range("D21").Value = '0.968'
$ws.Range("E21").Value = '  -2.92%  '

$ws.Range("D22").Value = '17.85'
$ws.Range("E22").Value = '  -0.46%  '

$ws.Range("D23").Value = '108.04'
$ws.Range("E23").Value = '  +4.98%  '

$ws.Range("D24").Value = '5.24'
$ws.Range("E24").Value = '  +1.68%  '

$ws.Range("D25").Value = '4.42'
$ws.Range("E25").Value = '  -4.97%  '

$ws.Range("E26").Value = '  -4.86%  '

$ws.Range("D27").Value = '5.99'
$ws.Range("E27").Value = '  -1.19%  '

$ws.Range("D28").Value = '10.34'
$ws.Range("E28").Value = '  -4.93%  '

$ws.Range("D29").Value = '9.72'
$ws.Range("E29").Value = '  +1.19%  '

$ws.Range("D30").Value = '33.41'
$ws.Range("E30").Value = '  -2.19%  '

$ws.Range("D31").Value = '4.39'
$ws.Range("E31").Value = '  +1.93%  '

$ws.Range("D32").Value = '6.91'
$ws.Range("E32").Value = '  -3.08%  '

$ws.Range("D33").Value = '11.93'
$ws.Range("E33").Value = '  -3.17%  '

$ws.Range("D34").Value = '0.112'
$ws.Range("E34").Value = '  -4.52%  '

$ws.Range("D35").Value = '62.30'
$ws.Range("E35").Value = '  -2.05%  '

$ws.Range("D36").Value = '3.807.04'
$ws.Range("E36").Value = '  -2.72%  '

$ws.Range("D37").Value = '0.0₃0820'
$ws.Range("E37").Value = '  -3.44%  '

$ws.Range("D38").Value = '1.00'
$ws.Range("E38").Value = '  -0.02%  '

$ws.Range("B39").Value = 'Bittensor'
$ws.Range("C39").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D39").Value = '508.40'
$ws.Range("E39").Value = '  -3.97%  '

$ws.Range("B40").Value = 'Stacks'
$ws.Range("C40").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D40").Value = '3.61'
$ws.Range("E40").Value = '  +1.00%  '

$ws.Range("D41").Value = '2.97'
$ws.Range("E41").Value = '  -6.43%  '

$ws.Range("B42").Value = 'Kaspa'
$ws.Range("C42").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D42").Value = '0.136'
$ws.Range("E42").Value = '  +1.07%  '

$ws.Range("B43").Value = 'TheGraph'
$ws.Range("C43").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D43").Value = '0.373'
$ws.Range("E43").Value = '  -5.32%  '

$ws.Range("D44").Value = '34.90'
$ws.Range("E44").Value = '  -6.19%  '

$ws.Range("D45").Value = '0.0458'
$ws.Range("E45").Value = '  -0.08%  '

$ws.Range("D46").Value = '2.96'
$ws.Range("E46").Value = '  +3.59%  '

$ws.Range("D47").Value = '3.39'
$ws.Range("E47").Value = '  +1.19%  '

$ws.Range("E48").Value = '  -2.61%  '

$ws.Range("E49").Value = '  -0.34%  '

$ws.Range("D50").Value = '8.36'
$ws.Range("E50").Value = '  -2.98%  '

$ws.Range("D51").Value = '1.81'
$ws.Range("E51").Value = '  +22.03%  '

# Reset style back to Normal (index 0) now that the text values are in place, so the
# temporary text number-format does not leave a stray style on the cells.
$priceRange.Style = "Normal"
